# Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics
# (percentages, dollar amounts) across the resume's achievements and
# work-experience bullet points, splitting the affected runs exactly as
# Word would when formatting a Find.Execute-selected sub-range.

$d = $word.ActiveDocument

function Set-MetricFormat($rng) {
    $rng.Font.Bold = 1
    $rng.Font.Color = 5258796   # RGB(0x2C, 0x3E, 0x50) -> 0x503E2C
}

# Highlights each token (in order of first appearance) within a paragraph's
# Range, scoping every Find.Execute call to "from here to the paragraph end"
# so repeated tokens elsewhere in the document are never touched.
function Highlight-Tokens($paraRange, [string[]]$tokens) {
    $cursor = $paraRange.Start
    foreach ($tok in $tokens) {
        $searchRange = $d.Range($cursor, $paraRange.End)
        $ok = $searchRange.Find.Execute($tok, $false, $false, $false, $false,
                                         $false, $true, 1, $false, "", 0)
        if ($ok) {
            Set-MetricFormat $searchRange
            $cursor = $searchRange.End
        }
    }
}

$PM = [char]0x00B1   # "±"

# "• Discovered systematic race coding errors ... from 23% to 64%"
Highlight-Tokens $d.Paragraphs.Item(10).Range @("23%", "64%")

# "• Utilized advanced sampling methods ... from ±4.2% to ±2.1%, increasing
#   voter turnout prediction accuracy from 71% to 87%, and ensuring ..."
Highlight-Tokens $d.Paragraphs.Item(12).Range @("${PM}4.2%", "${PM}2.1%", "71%", "87%")

# "• Trigonometric algorithm ... reduced mapping costs by 73.5%, saving
#   campaigns and organizations $4.7M and enabling ..."
Highlight-Tokens $d.Paragraphs.Item(13).Range @("73.5%", "$4.7M")

# "• Built real-time FEC analysis systems ... valued over $2 trillion"
Highlight-Tokens $d.Paragraphs.Item(14).Range @("$2")

# "• Modernized legacy ETL processes ... reducing processing time by 57%"
Highlight-Tokens $d.Paragraphs.Item(19).Range @("57%")

# "• Algorithmic innovation: Pioneered trigonometric boundary estimation
#   reducing mapping costs 73.5%"
Highlight-Tokens $d.Paragraphs.Item(55).Range @("73.5%")

# "• $4.7M savings enabled nonprofit access"
Highlight-Tokens $d.Paragraphs.Item(56).Range @("$4.7M")

# "• 178% accuracy improvement in racial classification algorithms"
Highlight-Tokens $d.Paragraphs.Item(58).Range @("178%")
